# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> bound to the slide master ("Integral" colours)
#   ppt/theme/theme2.xml  -> bound to the notes master  ("Office Theme" colours)
#
# The commit swaps the two themes' content: the slide master becomes the
# stock "Office Theme" colour scheme (dk1/lt1/dk2/lt2/accent1-6/hlink/
# folHlink), while the notes master ends up with the colours that used to
# belong to "Integral". Font scheme and format scheme were already
# byte-identical between the two theme parts, so only the 12-slot colour
# scheme actually changes.
#
# PowerPoint's object model exposes the presentation's theme colour scheme
# through Slide.ThemeColorScheme (a 12-item ColorScheme: dk1, lt1, dk2,
# lt2, accent1-6, hlink, folHlink, in that order) — apply the new "Office
# Theme" palette through it.

$p = $ppt.ActivePresentation

function ToOle([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB  = ToOle 0x00 0x00 0x00   # dk1      000000
$tcs.Item(2).RGB  = ToOle 0xFF 0xFF 0xFF   # lt1      FFFFFF
$tcs.Item(3).RGB  = ToOle 0x44 0x54 0x6A   # dk2      44546A
$tcs.Item(4).RGB  = ToOle 0xE7 0xE6 0xE6   # lt2      E7E6E6
$tcs.Item(5).RGB  = ToOle 0x5B 0x9B 0xD5   # accent1  5B9BD5
$tcs.Item(6).RGB  = ToOle 0xED 0x7D 0x31   # accent2  ED7D31
$tcs.Item(7).RGB  = ToOle 0xA5 0xA5 0xA5   # accent3  A5A5A5
$tcs.Item(8).RGB  = ToOle 0xFF 0xC0 0x00   # accent4  FFC000
$tcs.Item(9).RGB  = ToOle 0x44 0x72 0xC4   # accent5  4472C4
$tcs.Item(10).RGB = ToOle 0x70 0xAD 0x47   # accent6  70AD47
$tcs.Item(11).RGB = ToOle 0x05 0x63 0xC1   # hlink    0563C1
$tcs.Item(12).RGB = ToOle 0x95 0x4F 0x72   # folHlink 954F72
